$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "529.10"
# or "72.041.25" are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Row swaps (coin reordering) ---
# Rows 32/33: Hedera <-> Bittensor
$ws.Cells.Item(32, 2).Value = "Bittensor"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(32, 4).Value = "679.30"
$ws.Cells.Item(32, 5).Value = "  -2.91%  "

$ws.Cells.Item(33, 2).Value = "Hedera"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(33, 4).Value = "0.129"
$ws.Cells.Item(33, 5).Value = "  -2.11%  "

# Rows 46/47: Fetch.AI <-> ApeXProtocol
$ws.Cells.Item(46, 2).Value = "ApeXProtocol"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(46, 4).Value = "3.43"
$ws.Cells.Item(46, 5).Value = "  -2.84%  "

$ws.Cells.Item(47, 2).Value = "Fetch.AI"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(47, 4).Value = "2.63"
$ws.Cells.Item(47, 5).Value = "  -7.84%  "

# Rows 50/51: FLOKI <-> Monero
$ws.Cells.Item(50, 2).Value = "Monero"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(50, 4).Value = "146.52"
$ws.Cells.Item(50, 5).Value = "  +2.07%  "

$ws.Cells.Item(51, 2).Value = "FLOKI"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Cells.Item(51, 4).Value = "0.000270"
$ws.Cells.Item(51, 5).Value = "  -1.94%  "

# --- Price (column D) and Volume(1h) (column E) updates ---
$ws.Cells.Item(2, 4).Value = "72.041.25"
$ws.Cells.Item(2, 5).Value = "  +0.49%  "

$ws.Cells.Item(3, 4).Value = "4.011.40"
$ws.Cells.Item(3, 5).Value = "  -0.50%  "

$ws.Cells.Item(4, 5).Value = "  -0.04%  "

$ws.Cells.Item(5, 4).Value = "529.10"
$ws.Cells.Item(5, 5).Value = "  +0.71%  "

$ws.Cells.Item(6, 4).Value = "151.31"
$ws.Cells.Item(6, 5).Value = "  +1.91%  "

$ws.Cells.Item(7, 4).Value = "0.692"
$ws.Cells.Item(7, 5).Value = "  +10.49%  "

$ws.Cells.Item(8, 5).Value = "  +0.02%  "

$ws.Cells.Item(9, 4).Value = "0.746"
$ws.Cells.Item(9, 5).Value = "  +0.77%  "

$ws.Cells.Item(10, 4).Value = "0.171"
$ws.Cells.Item(10, 5).Value = "  -3.53%  "

$ws.Cells.Item(11, 5).Value = "  -3.89%  "

$ws.Cells.Item(12, 4).Value = "47.60"
$ws.Cells.Item(12, 5).Value = "  +3.09%  "

$ws.Cells.Item(13, 4).Value = "10.63"
$ws.Cells.Item(13, 5).Value = "  -1.82%  "

$ws.Cells.Item(14, 4).Value = "4.648.07"
$ws.Cells.Item(14, 5).Value = "  -0.38%  "

$ws.Cells.Item(15, 4).Value = "4.003.05"
$ws.Cells.Item(15, 5).Value = "  -0.83%  "

$ws.Cells.Item(16, 4).Value = "14.05"
$ws.Cells.Item(16, 5).Value = "  -1.96%  "

$ws.Cells.Item(17, 4).Value = "20.59"
$ws.Cells.Item(17, 5).Value = "  -4.04%  "

$ws.Cells.Item(18, 5).Value = "  -0.95%  "

$ws.Cells.Item(19, 5).Value = "  -3.26%  "

$ws.Cells.Item(20, 4).Value = "71.801.48"
$ws.Cells.Item(20, 5).Value = "  +0.21%  "

$ws.Cells.Item(21, 4).Value = "427.31"
$ws.Cells.Item(21, 5).Value = "  -3.62%  "

$ws.Cells.Item(22, 4).Value = "97.76"
$ws.Cells.Item(22, 5).Value = "  +2.67%  "

$ws.Cells.Item(23, 5).Value = "  -3.20%  "

$ws.Cells.Item(24, 5).Value = "  +3.26%  "

$ws.Cells.Item(25, 4).Value = "14.35"
$ws.Cells.Item(25, 5).Value = "  -0.34%  "

$ws.Cells.Item(26, 4).Value = "11.26"
$ws.Cells.Item(26, 5).Value = "  -8.45%  "

$ws.Cells.Item(27, 4).Value = "10.72"
$ws.Cells.Item(27, 5).Value = "  -3.44%  "

$ws.Cells.Item(28, 4).Value = "5.83"
$ws.Cells.Item(28, 5).Value = "  +0.76%  "

$ws.Cells.Item(29, 4).Value = "36.60"
$ws.Cells.Item(29, 5).Value = "  -1.51%  "

$ws.Cells.Item(30, 4).Value = "3.57"
$ws.Cells.Item(30, 5).Value = "  +22.28%  "

$ws.Cells.Item(31, 4).Value = "13.34"
$ws.Cells.Item(31, 5).Value = "  -2.52%  "

$ws.Cells.Item(34, 4).Value = "7.10"
$ws.Cells.Item(34, 5).Value = "  +1.92%  "

$ws.Cells.Item(35, 4).Value = "44.47"
$ws.Cells.Item(35, 5).Value = "  +7.49%  "

$ws.Cells.Item(36, 4).Value = "65.69"
$ws.Cells.Item(36, 5).Value = "  -3.32%  "

$ws.Cells.Item(37, 4).Value = "0.439"
$ws.Cells.Item(37, 5).Value = "  -2.19%  "

$ws.Cells.Item(38, 5).Value = "  -1.64%  "

$ws.Cells.Item(39, 4).Value = "0.0₃0827"
$ws.Cells.Item(39, 5).Value = "  -7.98%  "

$ws.Cells.Item(40, 5).Value = "  -3.79%  "

$ws.Cells.Item(41, 5).Value = "  -0.02%  "

$ws.Cells.Item(42, 5).Value = "  -0.10%  "

$ws.Cells.Item(43, 4).Value = "0.0487"
$ws.Cells.Item(43, 5).Value = "  -1.27%  "

$ws.Cells.Item(44, 4).Value = "3.20"
$ws.Cells.Item(44, 5).Value = "  +2.12%  "

$ws.Cells.Item(45, 5).Value = "  +2.12%  "

$ws.Cells.Item(48, 4).Value = "9.61"
$ws.Cells.Item(48, 5).Value = "  +3.32%  "

$ws.Cells.Item(49, 5).Value = "  -6.36%  "
